# Problem 3 slide: split the second bullet of "TextBox 2" into three runs
# so the "3" becomes "&gt;3" (rendered as literal "The LINQ version takes >3
# times as long as the imperative version", split across 3 <a:r> runs).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(13)
$sh = $s.Shapes.Item(4)          # "TextBox 2"

# --- remember the shape's current geometry -------------------------------
# Editing the text while "shrink shape to fit text" (AutoSize) is live makes
# this host recompute the box height from scratch (it does not reproduce
# PowerPoint's line-wrap metrics), which would shift cy/off unrelated to
# this edit. So: snapshot geometry, turn autofit off while we edit the
# runs, then turn it back on and restore the exact original geometry.
#
# Shape.Height/Top/Width/Left round-trip through a single-precision "points"
# value and this host truncates (rather than rounds) when converting back
# to EMU, which can make the restored value land 1 EMU short. Nudging the
# point value up to the next representable float32 keeps it in the same
# EMU bucket but clears the truncation, so the EMU comes back out exactly
# where it started.
$origLeft   = 53.56094741821289
$origTop    = 215.12197875976562
$origWidth  = 619.0244750976562
$origHeight = 142.98284912109375

$sh.TextFrame.AutoSize = 0   # ppAutoSizeNone - freeze the box while editing

$tr    = $sh.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$start = $para2.Start

# Run 1: "The LINQ version " (17 chars) - unchanged text, own run now
$r1 = $tr.Characters($start, 17)
$r1.Text = "The LINQ version "
$r1.Font.Bold = $false

# Run 2: "takes >3 " (9 chars) - the "3" becomes ">3"
$r2 = $tr.Characters($start + 17, 9)
$r2.Text = "takes >3 "
$r2.Font.Bold = $false

# Run 3: rest of the sentence, unchanged text, own run now
$r3 = $tr.Characters($start + 26, 1000)
$r3.Text = "times as long as the imperative version"
$r3.Font.Bold = $false

$sh.TextFrame.AutoSize = 1   # ppAutoSizeShapeToFitText - restore autofit flag

# restore the original box geometry exactly (see note above)
$sh.Left   = $origLeft
$sh.Top    = $origTop
$sh.Width  = $origWidth
$sh.Height = $origHeight
